# Applies the 2024-09-27 cryptos-list refresh:
#  - updates Price (D) / Volume 1h % (E) figures for most rows
#  - re-orders two pairs of rows whose ranking swapped
#    (Stacks/FirstDigitalUSD at 40/41, Aave/OKB at 43/44)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.546.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.68%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.649.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.24%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.649.63'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  +7.52%  '

$ws.Range("E11").Value = '  +2.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.87'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.62%  '

$ws.Range("E13").Value = '  +1.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000196'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +14.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.127.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.286.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.661.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.89'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '358.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.45%  '

$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.87%  '

$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.03%  '

$ws.Range("E27").Value = '  +15.12%  '

$ws.Range("E28").Value = '  -2.49%  '

$ws.Range("E29").Value = '  +2.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.22%  '

$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '527.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.57%  '

$ws.Range("E37").Value = '  +2.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '162.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.22%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.64%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.03%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '165.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.34%  '

$ws.Range("E47").Value = '  +2.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.17%  '

$ws.Range("E49").Value = '  +1.44%  '

$ws.Range("E50").Value = '  +3.60%  '

$ws.Range("E51").Value = '  +0.18%  '
